$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the ACTION columns D:E for the rule table (rows 18-25); the
# remaining ACTION column is consolidated into column C.
$ws.Range("D18:E25").ClearContents()

# Row 18 (header row): NAME / CONDITION / ACTION
$ws.Range("A18").Value = "NAME"
$ws.Range("B18").Value = "CONDITION"
$ws.Range("C18").Value = "ACTION"

# Row 19
$ws.Range("A19").Value = ""
$ws.Range("B19").Value = "Code changed 10010018"
$ws.Range("C19").Value = "Test"

# Row 20
$ws.Range("A20").Value = "New Teen Rule"
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = "gfndnvbx"

# Row 21 (B21 used to be a boolean FALSE; now a blank text cell)
$ws.Range("A21").Value = "Adult non-member"
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = "dsgagass"

# Row 22 (B22 used to be a boolean TRUE; now a blank text cell)
$ws.Range("A22").Value = "Adult member"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "dsvsbsb"

# Row 23
$ws.Range("A23").Value = "Senior perk"
$ws.Range("B23").Value = ""
$ws.Range("C23").Value = "nsngnsg"

# Row 24
$ws.Range("A24").Value = "Code changed 10010018"
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = "Code changed 10010018"

# Row 25
$ws.Range("A25").Value = "Senior perk"
$ws.Range("B25").Value = "Senior perk"
$ws.Range("C25").ClearContents()
